$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 100.16
$ws.Range("I33").Value = 100.16
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 100.16
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 128.84

$ws.Range("H107").Value = 1149.8
$ws.Range("I107").Value = 437.25
$ws.Range("J107").Value = 4000
$ws.Range("K107").Value = 437.25
$ws.Range("L107").Value = 4000
$ws.Range("M107").Value = 1482.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 452.625
$ws.Range("I32").Value = 452.625
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 452.625
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -165.625

$ws.Range("H45").Value = 1782.2727
$ws.Range("I45").Value = 1782.2727
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1782.2727
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1405.2727

$ws.Range("H61").Value = 1570.3334
$ws.Range("I61").Value = 1570.3334
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1570.3334
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1358.3334

$ws.Range("H74").Value = 3644.8
$ws.Range("I74").Value = 3644.8
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 3644.8
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -2770.8

$ws.Range("H77").Value = 3644.8
$ws.Range("I77").Value = 3644.8
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 18224
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -13856

$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws.Range("H102").Value = 4509.5
$ws.Range("I102").Value = 4509.5
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 4509.5
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -2887.5

$ws.Range("H104").Value = 58750
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 58750
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 58750
$ws.Range("N104").Value = -65738

$ws.Range("H106").Value = 27999
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 27999
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 27999
$ws.Range("N106").Value = -30523

$ws.Range("H110").Value = 582
$ws.Range("I110").Value = 632.6667
$ws.Range("J110").Value = 506
$ws.Range("K110").Value = 632.6667
$ws.Range("L110").Value = 506
$ws.Range("M110").Value = 1412.3333
$ws.Range("N110").Value = -4596

$ws.Range("H135").Value = 10000000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 10000000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 10000000
$ws.Range("N135").Value = -10010140

$ws.Range("H136").Value = 1570.3334
$ws.Range("I136").Value = 1570.3334
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4711.0002
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2161.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 9333
$ws.Range("I19").Value = 8000
$ws.Range("J19").Value = 9999.5
$ws.Range("K19").Value = 8000
$ws.Range("L19").Value = 9999.5
$ws.Range("M19").Value = -7827
$ws.Range("N19").Value = -10345.5

$ws.Range("H99").Value = 7999.5
$ws.Range("I99").Value = 7999.5
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 7999.5
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -6501.5

$ws.Range("H105").Value = 1999
$ws.Range("I105").Value = 1999
$ws.Range("J105").Value = 1999
$ws.Range("K105").Value = 1999
$ws.Range("L105").Value = 1999
$ws.Range("M105").Value = -252

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 5000
$ws.Range("I62").Value = 5000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 5000
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -4376
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 5000
$ws.Range("I65").Value = 5000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 25000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -21880
$ws.Range("N65").ClearContents()

$ws.Range("H122").Value = 1300
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1300
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 3900
$ws.Range("N122").Value = -8800

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 591.3333
$ws.Range("I18").Value = 704.75
$ws.Range("J18").Value = 364.5
$ws.Range("K18").Value = 2114.25
$ws.Range("L18").Value = 1093.5
$ws.Range("M18").Value = -1945.25
$ws.Range("N18").Value = -1431.5

$ws.Range("H86").Value = 2179.3
$ws.Range("I86").Value = 2065.889
$ws.Range("J86").Value = 3200
$ws.Range("K86").Value = 6197.667
$ws.Range("L86").Value = 9600
$ws.Range("M86").Value = -5011.667

$ws.Range("H89").Value = 2179.3
$ws.Range("I89").Value = 2065.889
$ws.Range("J89").Value = 3200
$ws.Range("K89").Value = 18593.001
$ws.Range("L89").Value = 28800
$ws.Range("M89").Value = -12665.001

$ws.Range("H98").Value = 1987.9
$ws.Range("I98").Value = 1822.7142
$ws.Range("J98").Value = 2373.3333
$ws.Range("K98").Value = 5468.142599999999
$ws.Range("L98").Value = 7119.999899999999
$ws.Range("M98").Value = -3970.142599999999
$ws.Range("N98").Value = -10115.9999

$ws.Range("H121").Value = 1380.6
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 1380.6
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 4141.799999999999
$ws.Range("N121").Value = -6761.799999999999
$ws.Range("M121").ClearContents()

$ws.Range("H131").Value = 2553.625
$ws.Range("I131").Value = 2205.8
$ws.Range("J131").Value = 3133.3333
$ws.Range("K131").Value = 6617.400000000001
$ws.Range("L131").Value = 9399.999899999999
$ws.Range("M131").Value = -1577.400000000001
$ws.Range("N131").Value = -19479.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1441.8572
$ws.Range("I97").Value = 1398.8
$ws.Range("J97").Value = 1549.5
$ws.Range("K97").Value = 1398.8
$ws.Range("L97").Value = 1549.5
$ws.Range("M97").Value = -902.8
$ws.Range("N97").Value = -2541.5

$ws.Range("H101").Value = 30000
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 30000
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 30000
$ws.Range("N101").Value = -36490

$ws.Range("H113").Value = 1900
$ws.Range("I113").Value = 1800
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1800
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 370

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 402199.2
$ws.Range("I46").Value = 1000750
$ws.Range("J46").Value = 3165.3333
$ws.Range("K46").Value = 1000750
$ws.Range("L46").Value = 3165.3333
$ws.Range("M46").Value = -1000562
$ws.Range("N46").Value = -3541.3333

$ws.Range("H74").Value = 50000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 50000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 50000
$ws.Range("N74").Value = -51996

$ws.Range("H77").Value = 50000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 50000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 150000
$ws.Range("N77").Value = -159984

$ws.Range("H82").Value = 1100.1428
$ws.Range("I82").Value = 1060.2
$ws.Range("J82").Value = 1200
$ws.Range("K82").Value = 1060.2
$ws.Range("L82").Value = 1200
$ws.Range("M82").Value = -699.2
$ws.Range("N82").Value = -1922

$ws.Range("H85").Value = 1100.1428
$ws.Range("I85").Value = 1060.2
$ws.Range("J85").Value = 1200
$ws.Range("K85").Value = 1060.2
$ws.Range("L85").Value = 1200
$ws.Range("M85").Value = 187.8
$ws.Range("N85").Value = -3696

$ws.Range("H106").Value = 9999
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 9999
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 9999
$ws.Range("N106").Value = -12523

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3000
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2376

$ws.Range("H65").Value = 3000
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -11880

$ws.Range("H75").Value = 50000
$ws.Range("I75").Value = 50000
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 50000
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -49064

$ws.Range("H78").Value = 50000
$ws.Range("I78").Value = 50000
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 150000
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -145320

$ws.Range("H117").Value = 21348.5
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 21348.5
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 21348.5
$ws.Range("N117").Value = -30526.5

$ws.Range("H122").Value = 3333.875
$ws.Range("I122").Value = 2976.4
$ws.Range("J122").Value = 3929.6667
$ws.Range("K122").Value = 8929.200000000001
$ws.Range("L122").Value = 11789.0001
$ws.Range("M122").Value = -6479.200000000001
